$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D15: change "@String "capacity change to newCapacity"" -> "log "capacity changed to newCapacity""
$ws.Range("D15").Value = "log “capacity changed to newCapacity”"

# Update D18: prefix "Log / @Error " instead of "@String "
$ws.Range("D18").Value = "Log / @Error “Plane ___ landed successfully.” / “Airport is full, cannot land.” / “This plane is already at this airport, cannot re-land.”"

# Update D19: prefix "Log / @Error " instead of "@String "
$ws.Range("D19").Value = "Log / @Error “Plane ___ has taken off from airport ____.” / “This plane is not at this airport, cannot take off.”"

# Update the active cell selection from D15 to D19
$ws.Range("D19").Select()
